$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) ---------------------------------------------------
# Entered in this specific order so new shared-string entries land at the
# same indices Excel produced (G2, E2, C2, H2, F2, D2).
$ws.Range("G2").Value = "GPU 1 Node(12 Core)"
$ws.Range("E2").Value = "CPU 1 Node(16 Core)"
$ws.Range("C2").Value = "Sparseness"
$ws.Range("H2").Value = "GPU 12 Node (+Head)"
$ws.Range("F2").Value = "CPU 12 Node(16 Core)"
$ws.Range("D2").Value = "REF(RDD) 1 GPU"

# --- Data rows (C:H), row by row ------------------------------------------
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 535
$ws.Range("E3").Value = 1745
$ws.Range("G3").Value = 4809
$ws.Range("H3").Value = 4972

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 248
$ws.Range("E4").Value = 1003
$ws.Range("G4").Value = 2414
$ws.Range("H4").Value = 2420

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 216
$ws.Range("E5").Value = 42042
$ws.Range("G5").Value = 2213
$ws.Range("H5").Value = 2224

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 301
$ws.Range("E6").Value = 41615
$ws.Range("G6").Value = 2330
$ws.Range("H6").Value = 2235

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 318
$ws.Range("E7").Value = 44620
$ws.Range("G7").Value = 2276
$ws.Range("H7").Value = 2180

$ws.Range("C8").Value = 0.1
$ws.Range("D8").Value = 373
$ws.Range("E8").Value = 45099
$ws.Range("G8").Value = 2952
$ws.Range("H8").Value = 2694

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 6396
$ws.Range("E9").Value = 101450
$ws.Range("G9").Value = 3634
$ws.Range("H9").Value = 3736

$ws.Range("C10").Value = 0.1
$ws.Range("D10").Value = 6631
$ws.Range("E10").Value = 107202
$ws.Range("G10").Value = 6958
$ws.Range("H10").Value = 6866

$ws.Range("C11").Value = 0.1
$ws.Range("D11").Value = 7470
$ws.Range("E11").Value = 110310
$ws.Range("G11").Value = 37418
$ws.Range("H11").Value = 38125

$ws.Range("C12").Value = 0.01
$ws.Range("D12").Value = 17014
$ws.Range("E12").Value = 177254
$ws.Range("G12").Value = 654384
$ws.Range("H12").Value = 645526

$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 6617119
$ws.Range("G13").Value = 97630
$ws.Range("H13").Value = 15388

$ws.Range("C14").Value = 0.1
$ws.Range("E14").Value = 6651854
$ws.Range("G14").Value = 398520
$ws.Range("H14").Value = 708049

$ws.Range("C15").Value = 0.01
$ws.Range("E15").Value = 6772861

$ws.Range("C16").Value = 0.0001
$ws.Range("E16").Value = 10812616

$ws.Range("C17").Value = 0.0001
$ws.Range("E17").Value = 47242268

# --- Column widths for the two new narrow columns --------------------------
$ws.Range("C1:D1").ColumnWidth = 11.5546875

# --- Selection --------------------------------------------------------------
[void]$ws.Range("D3").Select()
